# Update the "想去人数" (interested count) figures on the "展览" and "全部类型"
# sheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 588
$ws1.Range("F4").Value = 433
$ws1.Range("F7").Value = 2511
$ws1.Range("F9").Value = 6700
$ws1.Range("F10").Value = 178

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 588
$ws4.Range("F4").Value = 433
$ws4.Range("F9").Value = 2511
$ws4.Range("F11").Value = 6700
$ws4.Range("F12").Value = 178
